$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "https://www.youtube.com/watch?v=0O2Rq4HJBxw"
$errMsg = "Error: Summarizing and translating - Error in step 'Summarizing and translating': [red]" + [char]0x274C + " Expressiveness translation of block 80 failed after 3 retries. Please check ``output/gpt_log/error.json`` for more details.[/red]"
$ws.Range("E2").Value = $errMsg

$ws.Range("A3").Value = "https://www.youtube.com/watch?v=Qw4l1w0rkjs"

$ws.Range("A4").Value = "https://www.youtube.com/watch?v=MGyygiXMzRk"
$ws.Range("E4").Value = "Done"
